# YTI-696: Correcting integration test data according to new implementation.
#
# The Members_test_2 sheet's "CODEEXTENSION" (column G) values need to be
# updated from the bare code reference ("testcodeNN") to the prefixed
# reference format ("code:testcodeNN"), reflecting the new implementation.
#
# Shared strings are introduced in the order the cells below are written,
# so the write order here matches the order they were authored in (odd
# rows first, then even rows, then the last row) to reproduce the same
# shared-string table layout as the source edit.

$wb = $excel.ActiveWorkbook

$membersTest2 = $wb.Worksheets.Item("Members_test_2")

$membersTest2.Range("G3").Value = "code:testcode01"
$membersTest2.Range("G5").Value = "code:testcode03"
$membersTest2.Range("G7").Value = "code:testcode05"
$membersTest2.Range("G9").Value = "code:testcode07"

$membersTest2.Range("G4").Value = "code:testcode02"
$membersTest2.Range("G6").Value = "code:testcode04"
$membersTest2.Range("G8").Value = "code:testcode06"
$membersTest2.Range("G10").Value = "code:testcode08"

$membersTest2.Range("G11").Value = "code:testcode09"

# Restore/update the per-sheet cell selections to match the saved view
# state of the edited workbook.
$codeSchemes = $wb.Worksheets.Item("CodeSchemes")
$codeSchemes.Range("K42").Select()

$membersTest1 = $wb.Worksheets.Item("Members_test_1")
$membersTest1.Range("E2").Select()

# Members_test_2 is the last sheet interacted with / saved as the active
# tab, so select it (and a cell on it) last.
$membersTest2.Range("G12").Select()
